# Add status to demonstrate DEFAULT VALUE feature
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#Projects")

# Give the "Release board", "Odysseus" and "Argos" projects a status of
# "Completed" (rows 5-7 of the #Projects table, column C).
$ws.Range("C5").Value = "Completed"
$ws.Range("C6").Value = "Completed"
$ws.Range("C7").Value = "Completed"

# Make #Projects the active sheet again, with C8 selected.
$ws.Activate()
$ws.Range("C8").Select()
